# Temporarily disabling feedback loops until issues are sorted out.
$wb = $excel.ActiveWorkbook

# The BDMFL sheet holds the boolean lever in B2; flip it from 0 to 1 to
# disable the macroeconomic feedback loops.
$ws = $wb.Worksheets.Item("BDMFL")
$ws.Range("B2").Value = 1

# Make the BDMFL sheet the active/selected sheet, with B3 selected,
# matching the saved view state in the workbook.
$ws.Activate()
$ws.Range("B3").Select()
